$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.558.86'
$ws.Range('E2').Value = '  -2.36%  '
$ws.Range('D3').Value = '1.792.20'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.84'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5892'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2773'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06749'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.18'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07540'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.15%  '
$ws.Range('D12').Value = '1.798.54'
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.792'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6138'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.35%  '
$ws.Range('D15').Value = '2.035.90'
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '75.32'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008878'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.20%  '
$ws.Range('D18').Value = '28.549.94'
$ws.Range('E18').Value = '  -2.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.391'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '208.70'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.48'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.834'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.192'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1260'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.416'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06210'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.419'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.786'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.795'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.737'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.045'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6394'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.502'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.711'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01698'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.335'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('D41').Value = '1.142.02'
$ws.Range('E41').Value = '  -5.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8745'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.10%  '
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.16'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '1.944.99'
$ws.Range('E45').Value = '  -2.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000111'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.583'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.384'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05464'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4480'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.70%  '
